$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove bold styling previously applied to A1 (header 'Domain')
$ws.Range("A1").Font.Bold = $false

$data = @(
  @('Domain', 'Keywords', 'Title', 'Output'),
  @('http://cloudhiking.com', 'guides', 'cloudhiking - online trail guides and maps', '[cloudhiking - online trail guides and maps] (http://cloudhiking.com)'),
  @('http://finchglowholidays.com', 'destination,destinations', 'Finchglow Holidays – Book flights, hotels & accommodation, visit top destinations etc. with us', '[Finchglow Holidays – Book flights, hotels & accommodation, visit top destinations etc. with us] (http://finchglowholidays.com)'),
  @('http://bristol-hotels.net', 'accommodations', 'Bristol hotels & apartments, all accommodations in Bristol', '[Bristol hotels & apartments, all accommodations in Bristol] (http://bristol-hotels.net)'),
  @('http://eddingtonhouseinn.com', 'accommodations', 'Eddington House Inn - Bennington Vermont Bed and Breakfast, Bennington College lodging accommodations, North Bennington Village, Vermont Spa Getaway,VT Wine Tasting Dinner Packages, Bennington College Visits', '[Eddington House Inn - Bennington Vermont Bed and Breakfast, Bennington College lodging accommodations, North Bennington Village, Vermont Spa Getaway,VT Wine Tasting Dinner Packages, Bennington College Visits] (http://eddingtonhouseinn.com)'),
  @('http://floridatraveler.org', 'travel', 'floridatraveler | Take a trip to Florida: present, past, and future!', '[floridatraveler | Take a trip to Florida: present, past, and future!] (http://floridatraveler.org)'),
  @('http://cascadiakids.com', 'travel', 'Cascadia Kids : Family travel in the Pacific Northwest and BC', '[Cascadia Kids : Family travel in the Pacific Northwest and BC] (http://cascadiakids.com)')
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}
